$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 (weather refresh) ---
$ws.Cells.Item(2,2).Value = 23.09
$ws.Cells.Item(2,3).Value = 71.0
$ws.Cells.Item(2,4).Value = "few clouds"
$ws.Cells.Item(2,5).Value = "2025-07-18T22:43:02.840109"
$ws.Cells.Item(2,6).Value = "Berlin"

$ws.Cells.Item(3,2).Value = 24.03
$ws.Cells.Item(3,3).Value = 73.0
$ws.Cells.Item(3,4).Value = "clear sky"
$ws.Cells.Item(3,5).Value = "2025-07-18T22:43:02.971062"
$ws.Cells.Item(3,6).Value = "Baku"

$ws.Cells.Item(4,2).Value = 23.09
$ws.Cells.Item(4,3).Value = 71.0
$ws.Cells.Item(4,4).Value = "few clouds"
$ws.Cells.Item(4,5).Value = "2025-07-18T22:43:19.383522"
$ws.Cells.Item(4,6).Value = "Berlin"

$ws.Cells.Item(5,2).Value = 24.03
$ws.Cells.Item(5,3).Value = 73.0
$ws.Cells.Item(5,4).Value = "clear sky"
$ws.Cells.Item(5,5).Value = "2025-07-18T22:43:19.492466"
$ws.Cells.Item(5,6).Value = "Baku"

# --- Append new rows 6-8 ---
$ws.Cells.Item(6,1).Value = 5.0
$ws.Cells.Item(6,2).Value = 22.95
$ws.Cells.Item(6,3).Value = 71.0
$ws.Cells.Item(6,4).Value = "few clouds"
$ws.Cells.Item(6,5).Value = "2025-07-18T22:47:45.439110"
$ws.Cells.Item(6,6).Value = "Berlin"

$ws.Cells.Item(7,1).Value = 6.0
$ws.Cells.Item(7,2).Value = 24.03
$ws.Cells.Item(7,3).Value = 73.0
$ws.Cells.Item(7,4).Value = "clear sky"
$ws.Cells.Item(7,5).Value = "2025-07-18T22:47:45.581425"
$ws.Cells.Item(7,6).Value = "Baku"

$ws.Cells.Item(8,1).Value = 7.0
$ws.Cells.Item(8,2).Value = 19.24
$ws.Cells.Item(8,3).Value = 58.0
$ws.Cells.Item(8,4).Value = "overcast clouds"
$ws.Cells.Item(8,5).Value = "2025-07-18T22:47:45.691338"
$ws.Cells.Item(8,6).Value = "Moscow"

# --- Resize columns to fit the new (longer) content ---
# "Description" (D) now also holds the longer "overcast clouds" and
# "City" (F) now also holds "Moscow" - both columns need to widen
# (closest attainable width to the bestFit values Apache POI would compute).
$ws.Columns.Item(4).ColumnWidth = 12.5
$ws.Columns.Item(6).ColumnWidth = 6.75
